# Auto-generated edit script: updates FFXIV leve profit tracker values
# per scheduled market-price refresh (see commit message).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1079.9375
$ws.Cells.Item(15, 9).Value = 1079.9375
$ws.Cells.Item(15, 11).Value = 3239.8125
$ws.Cells.Item(15, 13).Value = -3070.8125
$ws.Cells.Item(32, 8).Value = 1279.1666
$ws.Cells.Item(32, 9).Value = 1295.6666
$ws.Cells.Item(32, 10).Value = 1262.6666
$ws.Cells.Item(32, 11).Value = 1295.6666
$ws.Cells.Item(32, 12).Value = 1262.6666
$ws.Cells.Item(32, 13).Value = -969.6666
$ws.Cells.Item(32, 14).Value = -1914.6666
$ws.Cells.Item(70, 8).Value = 31241858
$ws.Cells.Item(70, 9).Value = 12385600
$ws.Cells.Item(70, 10).Value = 59526250
$ws.Cells.Item(70, 11).Value = 37156800
$ws.Cells.Item(70, 12).Value = 178578750
$ws.Cells.Item(70, 13).Value = -37156530
$ws.Cells.Item(70, 14).Value = -178579290
$ws.Cells.Item(73, 8).Value = 31241858
$ws.Cells.Item(73, 9).Value = 12385600
$ws.Cells.Item(73, 10).Value = 59526250
$ws.Cells.Item(73, 11).Value = 37156800
$ws.Cells.Item(73, 12).Value = 178578750
$ws.Cells.Item(73, 13).Value = -37155864
$ws.Cells.Item(73, 14).Value = -178580622
$ws.Cells.Item(80, 8).Value = 4049433.2
$ws.Cells.Item(80, 9).Value = 6584898
$ws.Cells.Item(80, 11).Value = 19754694
$ws.Cells.Item(80, 13).Value = -19753696
$ws.Cells.Item(83, 8).Value = 4049433.2
$ws.Cells.Item(83, 9).Value = 6584898
$ws.Cells.Item(83, 11).Value = 59264082
$ws.Cells.Item(83, 13).Value = -59259090
$ws.Cells.Item(98, 8).Value = 27029356
$ws.Cells.Item(98, 9).Value = 27029356
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 27029356
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = -27027858
$ws.Cells.Item(98, 14).ClearContents()  # N98 removed (was -8551)
$ws.Cells.Item(106, 8).Value = 58825784
$ws.Cells.Item(106, 9).Value = 71430210
$ws.Cells.Item(106, 11).Value = 71430210
$ws.Cells.Item(106, 13).Value = -71429579
$ws.Cells.Item(116, 8).Value = 17864986
$ws.Cells.Item(116, 9).Value = 50003600
$ws.Cells.Item(116, 11).Value = 50003600
$ws.Cells.Item(116, 13).Value = -50000158
$ws.Cells.Item(122, 8).Value = 27029356
$ws.Cells.Item(122, 9).Value = 27029356
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 81088068
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -81085618
$ws.Cells.Item(122, 14).ClearContents()  # N122 removed (was -21565)
$ws.Cells.Item(132, 8).Value = 1981.4412
$ws.Cells.Item(132, 9).Value = 2016.5
$ws.Cells.Item(132, 10).Value = 1817.8334
$ws.Cells.Item(132, 11).Value = 6049.5
$ws.Cells.Item(132, 12).Value = 5453.5002
$ws.Cells.Item(132, 13).Value = -3519.5
$ws.Cells.Item(132, 14).Value = -10513.5002
$ws.Cells.Item(135, 8).Value = 435883.53
$ws.Cells.Item(135, 9).Value = 501083.06
$ws.Cells.Item(135, 11).Value = 4509747.54
$ws.Cells.Item(135, 13).Value = -4507212.54
$ws.Cells.Item(138, 8).Value = 1566364.2
$ws.Cells.Item(138, 9).Value = 739.375
$ws.Cells.Item(138, 10).Value = 2505739.2
$ws.Cells.Item(138, 11).Value = 2218.125
$ws.Cells.Item(138, 12).Value = 7517217.600000001
$ws.Cells.Item(138, 13).Value = 2921.875
$ws.Cells.Item(138, 14).Value = -7527497.600000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 58827012
$ws.Cells.Item(2, 10).Value = 111116480
$ws.Cells.Item(2, 12).Value = 111116480
$ws.Cells.Item(2, 14).Value = -111116706
$ws.Cells.Item(32, 8).Value = 2239483
$ws.Cells.Item(32, 9).Value = 2609296.5
$ws.Cells.Item(32, 11).Value = 2609296.5
$ws.Cells.Item(32, 13).Value = -2609009.5
$ws.Cells.Item(61, 8).Value = 7066.7095
$ws.Cells.Item(61, 9).Value = 3008.2778
$ws.Cells.Item(61, 10).Value = 12686.077
$ws.Cells.Item(61, 11).Value = 3008.2778
$ws.Cells.Item(61, 12).Value = 12686.077
$ws.Cells.Item(61, 13).Value = -2796.2778
$ws.Cells.Item(61, 14).Value = -13110.077
$ws.Cells.Item(102, 8).Value = 1729.4
$ws.Cells.Item(102, 9).Value = 1515.5
$ws.Cells.Item(102, 11).Value = 1515.5
$ws.Cells.Item(102, 13).Value = 106.5
$ws.Cells.Item(116, 8).Value = 58827012
$ws.Cells.Item(116, 10).Value = 111116480
$ws.Cells.Item(116, 12).Value = 111116480
$ws.Cells.Item(116, 14).Value = -111121068
$ws.Cells.Item(122, 8).Value = 21285.428
$ws.Cells.Item(122, 9).Value = 23817.908
$ws.Cells.Item(122, 10).Value = 11999.667
$ws.Cells.Item(122, 11).Value = 71453.724
$ws.Cells.Item(122, 12).Value = 35999.001
$ws.Cells.Item(122, 13).Value = -69003.724
$ws.Cells.Item(122, 14).Value = -40899.001
$ws.Cells.Item(132, 8).Value = 3329.2034
$ws.Cells.Item(132, 9).Value = 1954.2291
$ws.Cells.Item(132, 11).Value = 5862.6873
$ws.Cells.Item(132, 13).Value = -3332.6873
$ws.Cells.Item(136, 8).Value = 7066.7095
$ws.Cells.Item(136, 9).Value = 3008.2778
$ws.Cells.Item(136, 10).Value = 12686.077
$ws.Cells.Item(136, 11).Value = 9024.8334
$ws.Cells.Item(136, 12).Value = 38058.231
$ws.Cells.Item(136, 13).Value = -6474.8334
$ws.Cells.Item(136, 14).Value = -43158.231

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 58827012
$ws.Cells.Item(3, 10).Value = 111116480
$ws.Cells.Item(3, 12).Value = 111116480
$ws.Cells.Item(3, 14).Value = -111116708

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5854392
$ws.Cells.Item(31, 9).Value = 2792.25
$ws.Cells.Item(31, 10).Value = 13344440
$ws.Cells.Item(31, 11).Value = 2792.25
$ws.Cells.Item(31, 12).Value = 13344440
$ws.Cells.Item(31, 13).Value = -2497.25
$ws.Cells.Item(31, 14).Value = -13345030
$ws.Cells.Item(34, 8).Value = 5854392
$ws.Cells.Item(34, 9).Value = 2792.25
$ws.Cells.Item(34, 10).Value = 13344440
$ws.Cells.Item(34, 11).Value = 2792.25
$ws.Cells.Item(34, 12).Value = 13344440
$ws.Cells.Item(34, 13).Value = -2590.25
$ws.Cells.Item(34, 14).Value = -13344844
$ws.Cells.Item(58, 8).Value = 10644505
$ws.Cells.Item(58, 9).Value = 20002620
$ws.Cells.Item(58, 11).Value = 20002620
$ws.Cells.Item(58, 13).Value = -20002417
$ws.Cells.Item(86, 8).Value = 6949166
$ws.Cells.Item(86, 9).Value = 20834498
$ws.Cells.Item(86, 11).Value = 20834498
$ws.Cells.Item(86, 13).Value = -20833375
$ws.Cells.Item(89, 8).Value = 6949166
$ws.Cells.Item(89, 9).Value = 20834498
$ws.Cells.Item(89, 11).Value = 104172490
$ws.Cells.Item(89, 13).Value = -104166874
$ws.Cells.Item(99, 8).Value = 8078.8423
$ws.Cells.Item(99, 9).Value = 9942.429
$ws.Cells.Item(99, 11).Value = 9942.429
$ws.Cells.Item(99, 13).Value = -8444.429
$ws.Cells.Item(107, 8).Value = 2500.913
$ws.Cells.Item(107, 10).Value = 2316.3333
$ws.Cells.Item(107, 12).Value = 2316.3333
$ws.Cells.Item(107, 14).Value = -6156.3333
$ws.Cells.Item(122, 8).Value = 2210.35
$ws.Cells.Item(122, 9).Value = 1525.625
$ws.Cells.Item(122, 10).Value = 2666.8333
$ws.Cells.Item(122, 11).Value = 4576.875
$ws.Cells.Item(122, 12).Value = 8000.499899999999
$ws.Cells.Item(122, 13).Value = -2126.875
$ws.Cells.Item(122, 14).Value = -12900.4999
$ws.Cells.Item(126, 8).Value = 8078.8423
$ws.Cells.Item(126, 9).Value = 9942.429
$ws.Cells.Item(126, 11).Value = 29827.287
$ws.Cells.Item(126, 13).Value = -27357.287
$ws.Cells.Item(134, 8).Value = 5978.12
$ws.Cells.Item(134, 9).Value = 2822.6191
$ws.Cells.Item(134, 10).Value = 8263.138000000001
$ws.Cells.Item(134, 11).Value = 8467.8573
$ws.Cells.Item(134, 12).Value = 24789.414
$ws.Cells.Item(134, 13).Value = -5932.8573
$ws.Cells.Item(134, 14).Value = -29859.414
$ws.Cells.Item(136, 8).Value = 10644505
$ws.Cells.Item(136, 9).Value = 20002620
$ws.Cells.Item(136, 11).Value = 60007860
$ws.Cells.Item(136, 13).Value = -60005310

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(112, 8).Value = 4399.5
$ws.Cells.Item(112, 9).Value = 2999
$ws.Cells.Item(112, 11).Value = 8997
$ws.Cells.Item(112, 13).Value = -7889
$ws.Cells.Item(113, 8).Value = 2520.4375
$ws.Cells.Item(113, 9).Value = 1551
$ws.Cells.Item(113, 10).Value = 3102.1
$ws.Cells.Item(113, 11).Value = 4653
$ws.Cells.Item(113, 12).Value = 9306.299999999999
$ws.Cells.Item(113, 13).Value = -2483
$ws.Cells.Item(113, 14).Value = -13646.3
$ws.Cells.Item(132, 8).Value = 10810.821
$ws.Cells.Item(132, 9).Value = 4164.706
$ws.Cells.Item(132, 11).Value = 37482.354
$ws.Cells.Item(132, 13).Value = -34952.354
$ws.Cells.Item(140, 8).Value = 3744.5881
$ws.Cells.Item(140, 9).Value = 1708.5
$ws.Cells.Item(140, 10).Value = 5554.4443
$ws.Cells.Item(140, 11).Value = 5125.5
$ws.Cells.Item(140, 12).Value = 16663.3329
$ws.Cells.Item(140, 13).Value = 54.5
$ws.Cells.Item(140, 14).Value = -27023.3329
$ws.Cells.Item(141, 8).Value = 5585.72
$ws.Cells.Item(141, 9).Value = 2482.15
$ws.Cells.Item(141, 11).Value = 7446.450000000001
$ws.Cells.Item(141, 13).Value = -2266.450000000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(40, 8).Value = 21333.334
$ws.Cells.Item(43, 8).Value = 1500
$ws.Cells.Item(43, 9).Value = 1500
$ws.Cells.Item(43, 11).Value = 1500
$ws.Cells.Item(43, 13).Value = -1349
$ws.Cells.Item(124, 8).Value = 75540.5
$ws.Cells.Item(124, 10).Value = 75540.5
$ws.Cells.Item(124, 12).Value = 75540.5
$ws.Cells.Item(124, 14).Value = -85360.5
$ws.Cells.Item(132, 8).Value = 5184.436
$ws.Cells.Item(132, 9).Value = 2167.84
$ws.Cells.Item(132, 11).Value = 6503.52
$ws.Cells.Item(132, 13).Value = -3973.52

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 8999.625
$ws.Cells.Item(93, 8).Value = 5416.2256
$ws.Cells.Item(93, 10).Value = 8248.5
$ws.Cells.Item(93, 12).Value = 8248.5
$ws.Cells.Item(93, 14).Value = -10744.5
$ws.Cells.Item(121, 8).Value = 44582.4
$ws.Cells.Item(121, 10).Value = 44582.4
$ws.Cells.Item(121, 12).Value = 44582.4
$ws.Cells.Item(121, 14).Value = -48076.4
$ws.Cells.Item(126, 8).Value = 8999.625
$ws.Cells.Item(130, 8).Value = 78996
$ws.Cells.Item(130, 10).Value = 78996
$ws.Cells.Item(130, 12).Value = 78996
$ws.Cells.Item(130, 14).Value = -89036
$ws.Cells.Item(132, 8).Value = 7942216.5
$ws.Cells.Item(132, 9).Value = 15627534
$ws.Cells.Item(132, 11).Value = 46882602
$ws.Cells.Item(132, 13).Value = -46880072

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 13).ClearContents()  # M42 removed (was 228)
$ws.Cells.Item(122, 8).Value = 120480.414
$ws.Cells.Item(122, 9).Value = 168726.95
$ws.Cells.Item(122, 10).Value = 4688.7
$ws.Cells.Item(122, 11).Value = 506180.85
$ws.Cells.Item(122, 12).Value = 14066.1
$ws.Cells.Item(122, 13).Value = -503730.85
$ws.Cells.Item(122, 14).Value = -18966.1
$ws.Cells.Item(132, 8).Value = 83395660
$ws.Cells.Item(132, 9).Value = 250050000
$ws.Cells.Item(132, 10).Value = 68500
$ws.Cells.Item(132, 11).Value = 750150000
$ws.Cells.Item(132, 12).Value = 205500
$ws.Cells.Item(132, 13).Value = -750147470
$ws.Cells.Item(132, 14).Value = -210560
$ws.Cells.Item(136, 8).Value = 22754034
$ws.Cells.Item(136, 9).Value = 43479092
$ws.Cells.Item(136, 10).Value = 55157.332
$ws.Cells.Item(136, 11).Value = 130437276
$ws.Cells.Item(136, 12).Value = 165471.996
$ws.Cells.Item(136, 13).Value = -130434726
$ws.Cells.Item(136, 14).Value = -170571.996

